# 20250411 LC3WP + LC1EP
$wb = $excel.ActiveWorkbook

# Update the ISIN for ROCKWOOL B on the "Europe" sheet (row 15)
$wsEurope = $wb.Worksheets.Item("Europe")
$wsEurope.Range("B15").Value = "DK0063855168"

# Update the ISIN for ROCKWOOL B on the "World" sheet (row 15)
$wsWorld = $wb.Worksheets.Item("World")
$wsWorld.Range("B15").Value = "DK0063855168"

# Update cursor/selection on "World" sheet
$wsWorld.Activate() | Out-Null
$wsWorld.Range("E2:E14").Select() | Out-Null

# Make "Europe" the active sheet & set its new selection (was "World")
$wsEurope.Activate() | Out-Null
$wsEurope.Range("B15").Select() | Out-Null
